$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.467.24"
$ws.Range("E2").Value = "  +7.94%  "
$ws.Range("D3").Value = "3.409.29"
$ws.Range("E3").Value = "  +5.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'411.29"
$ws.Range("D6").Value = "'122.34"
$ws.Range("E6").Value = "  +13.40%  "
$ws.Range("D7").Value = "3.405.02"
$ws.Range("E7").Value = "  +5.05%  "
$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D10").Value = "'0.639"
$ws.Range("E10").Value = "  +3.55%  "
$ws.Range("D11").Value = "'0.120"
$ws.Range("E11").Value = "  +27.02%  "
$ws.Range("D12").Value = "'41.07"
$ws.Range("E12").Value = "  +4.98%  "
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "3.940.06"
$ws.Range("E14").Value = "  +4.95%  "
$ws.Range("D15").Value = "'8.40"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "'19.50"
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("D17").Value = "3.400.05"
$ws.Range("E17").Value = "  +4.96%  "
$ws.Range("D18").Value = "61.346.00"
$ws.Range("E18").Value = "  +8.05%  "
$ws.Range("D19").Value = "'1.02"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").Value = "'10.84"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  +12.08%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'12.81"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "'298.76"
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("D25").Value = "'76.48"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "'30.70"
$ws.Range("E27").Value = "  +9.73%  "
$ws.Range("E28").Value = "  +10.44%  "
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").Value = "'7.62"
$ws.Range("E30").Value = "  -6.21%  "
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'0.116"
$ws.Range("E32").Value = "  +5.71%  "
$ws.Range("D33").Value = "'42.63"
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.38"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "'2.53"
$ws.Range("E35").Value = "  +18.90%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "'0.0477"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "'52.33"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "'3.02"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("E42").Value = "  +5.49%  "
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("D44").Value = "'133.98"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "'17.13"
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.282"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.90"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "'2.18"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").Value = "'21.72"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "2.203.16"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").Value = "3.743.57"
$ws.Range("E51").Value = "  +5.02%  "
